# Insert a new data row at row 9 (pushing existing rows 9-21 down to 10-22)
# and populate it with the new weekly price entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(9).Insert()

$ws.Cells.Item(9, 1).Value = 1
$ws.Cells.Item(9, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(9, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(9, 4).Value = 44482
$ws.Cells.Item(9, 5).Value = 15
$ws.Cells.Item(9, 6).Value = 100112013
$ws.Cells.Item(9, 7).Value = "Alcachofa"
$ws.Cells.Item(9, 8).Value = "Madrigal"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 200
$ws.Cells.Item(9, 11).Value = 14000
$ws.Cells.Item(9, 12).Value = 15000
$ws.Cells.Item(9, 13).Value = 14500
$ws.Cells.Item(9, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(9, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(9, 16).Value = 362
$ws.Cells.Item(9, 17).Value = 40
$ws.Cells.Item(9, 18).Value = "Hortaliza"
